$d = $word.ActiveDocument

# --- Split the final paragraph into three paragraphs ------------------
# Before: last paragraph = "Will add feature ... entered" (ilvl=1) followed
#         by the hidden "_GoBack" bookmark.
# After : that paragraph keeps its text (bookmark removed from it), then a
#         new ilvl=1 paragraph "Too annoying to implement" and a new
#         ilvl=0 paragraph "Will work on decoder next" are appended, with
#         the "_GoBack" bookmark now sitting at the end of the new last
#         paragraph.

# 1. Append a new (empty) paragraph after the current last paragraph; it
#    inherits the ListParagraph / ilvl=1 / numId=1 formatting we want for
#    "Too annoying to implement".
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$paras = $d.Paragraphs
$annoyingPara = $paras.Last
$annoyingPara.Range.Text = "Too annoying to implement"

# 2. Append another new paragraph after that one for "Will work on decoder
#    next", promoted to ilvl=0 (ListLevelNumber is 1-based, so 1 == ilvl 0).
$annoyingPara.Range.InsertParagraphAfter()

$paras = $d.Paragraphs
$decoderPara = $paras.Last
$decoderPara.Range.ListFormat.ListLevelNumber = 1
# Type a trailing placeholder character so that the bookmark insertion
# point below is not the very last character before a paragraph mark
# (inserting a bookmark exactly there is mishandled by the host and snaps
# to the start of the document) -- we remove the placeholder afterwards.
$decoderPara.Range.Text = "Will work on decoder nextX"

# --- Move the "_GoBack" bookmark to the end of the new last paragraph -
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$paras = $d.Paragraphs
$decoderPara = $paras.Last
$placeholderPos = $decoderPara.Range.End - 2
$bookmarkTarget = $d.Range($placeholderPos, $placeholderPos)
$d.Bookmarks.Add("_GoBack", $bookmarkTarget)

# Remove the placeholder character now that the bookmark is anchored.
$placeholderRange = $d.Range($placeholderPos, $placeholderPos + 1)
$placeholderRange.Delete()
